$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1113.1666
$ws.Range("J40").Value = 1500
$ws.Range("L40").Value = 1500
$ws.Range("N40").Value = -1850

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2317.9644
$ws.Range("I74").Value = 1670.15
$ws.Range("J74").Value = 3937.5
$ws.Range("K74").Value = 1670.15
$ws.Range("L74").Value = 3937.5
$ws.Range("M74").Value = -734.1500000000001
$ws.Range("N74").Value = -5809.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 2317.9644
$ws.Range("I77").Value = 1670.15
$ws.Range("J77").Value = 3937.5
$ws.Range("K77").Value = 8350.75
$ws.Range("L77").Value = 19687.5
$ws.Range("M77").Value = -3670.75
$ws.Range("N77").Value = -29047.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2493.75
$ws.Range("I100").Value = 1860
$ws.Range("K100").Value = 1860
$ws.Range("M100").Value = -1319

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1525.5
$ws.Range("I106").Value = 1305.5555
$ws.Range("K106").Value = 1305.5555
$ws.Range("M106").Value = -674.5554999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4313.9165
$ws.Range("I113").Value = 3268.5715
$ws.Range("J113").Value = 5777.4
$ws.Range("K113").Value = 3268.5715
$ws.Range("L113").Value = 5777.4
$ws.Range("M113").Value = -14.57150000000001
$ws.Range("N113").Value = -12285.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2631.1
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 2701.2222
$ws.Range("K116").Value = 2000
$ws.Range("L116").Value = 2701.2222
$ws.Range("M116").Value = 1442
$ws.Range("N116").Value = -9585.2222

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 53761.75
$ws.Range("I132").Value = 59679.723
$ws.Range("K132").Value = 179039.169
$ws.Range("M132").Value = -176509.169

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 38508.89
$ws.Range("I137").Value = 1214.5385
$ws.Range("J137").Value = 73139.36
$ws.Range("K137").Value = 3643.6155
$ws.Range("L137").Value = 219418.08
$ws.Range("M137").Value = -1093.6155
$ws.Range("N137").Value = -224518.08

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2980
$ws.Range("I2").Value = 1711.75
$ws.Range("J2").Value = 4671
$ws.Range("K2").Value = 1711.75
$ws.Range("L2").Value = 4671
$ws.Range("M2").Value = -1598.75
$ws.Range("N2").Value = -4897

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20621.75
$ws.Range("I32").Value = 25009.488
$ws.Range("J32").Value = 6108.4614
$ws.Range("K32").Value = 25009.488
$ws.Range("L32").Value = 6108.4614
$ws.Range("M32").Value = -24722.488
$ws.Range("N32").Value = -6682.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2260.682
$ws.Range("I45").Value = 1876
$ws.Range("K45").Value = 1876
$ws.Range("M45").Value = -1499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 4499.7144
$ws.Range("I97").Value = 4603.8
$ws.Range("K97").Value = 4603.8
$ws.Range("M97").Value = -4107.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2980
$ws.Range("I116").Value = 1711.75
$ws.Range("J116").Value = 4671
$ws.Range("K116").Value = 1711.75
$ws.Range("L116").Value = 4671
$ws.Range("M116").Value = 582.25
$ws.Range("N116").Value = -9259

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 29049.475
$ws.Range("I132").Value = 2804.4443
$ws.Range("J132").Value = 52670
$ws.Range("K132").Value = 8413.332900000001
$ws.Range("L132").Value = 158010
$ws.Range("M132").Value = -5883.332900000001
$ws.Range("N132").Value = -163070

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2980
$ws.Range("I3").Value = 1711.75
$ws.Range("J3").Value = 4671
$ws.Range("K3").Value = 1711.75
$ws.Range("L3").Value = 4671
$ws.Range("M3").Value = -1597.75
$ws.Range("N3").Value = -4899

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12617.517
$ws.Range("I31").Value = 37266.777
$ws.Range("J31").Value = 2533.7273
$ws.Range("K31").Value = 37266.777
$ws.Range("L31").Value = 2533.7273
$ws.Range("M31").Value = -36971.777
$ws.Range("N31").Value = -3123.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 12617.517
$ws.Range("I34").Value = 37266.777
$ws.Range("J34").Value = 2533.7273
$ws.Range("K34").Value = 37266.777
$ws.Range("L34").Value = 2533.7273
$ws.Range("M34").Value = -37064.777
$ws.Range("N34").Value = -2937.7273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 584.04443
$ws.Range("I5").Value = 473.5
$ws.Range("J5").Value = 672.48
$ws.Range("K5").Value = 1420.5
$ws.Range("L5").Value = 2017.44
$ws.Range("M5").Value = -1308.5
$ws.Range("N5").Value = -2241.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 356.5238
$ws.Range("J122").Value = 370.83334
$ws.Range("L122").Value = 3337.50006
$ws.Range("N122").Value = -8237.50006

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 137829.39
$ws.Range("J131").Value = 147891.94
$ws.Range("L131").Value = 443675.82
$ws.Range("N131").Value = -453755.82

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 584.04443
$ws.Range("I135").Value = 473.5
$ws.Range("J135").Value = 672.48
$ws.Range("K135").Value = 4261.5
$ws.Range("L135").Value = 6052.32
$ws.Range("M135").Value = -1726.5
$ws.Range("N135").Value = -11122.32

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7844.2856
$ws.Range("I80").Value = 12215
$ws.Range("J80").Value = 3870.9092
$ws.Range("K80").Value = 12215
$ws.Range("L80").Value = 3870.9092
$ws.Range("M80").Value = -11217
$ws.Range("N80").Value = -5866.9092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 7844.2856
$ws.Range("I83").Value = 12215
$ws.Range("J83").Value = 3870.9092
$ws.Range("K83").Value = 61075
$ws.Range("L83").Value = 19354.546
$ws.Range("M83").Value = -56083
$ws.Range("N83").Value = -29338.546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2636.3635
$ws.Range("I102").Value = 2766.6667
$ws.Range("J102").Value = 2480
$ws.Range("K102").Value = 2766.6667
$ws.Range("L102").Value = 2480
$ws.Range("M102").Value = -1144.6667
$ws.Range("N102").Value = -5724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3164.72
$ws.Range("I113").Value = 2500.9
$ws.Range("K113").Value = 2500.9
$ws.Range("M113").Value = -330.9000000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5270.9688
$ws.Range("I126").Value = 5156.05
$ws.Range("J126").Value = 5462.5
$ws.Range("K126").Value = 15468.15
$ws.Range("L126").Value = 16387.5
$ws.Range("M126").Value = -12998.15
$ws.Range("N126").Value = -21327.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3619
$ws.Range("I122").Value = 3120.8
$ws.Range("J122").Value = 4034.1667
$ws.Range("K122").Value = 9362.400000000001
$ws.Range("L122").Value = 12102.5001
$ws.Range("M122").Value = -6912.400000000001
$ws.Range("N122").Value = -17002.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1273.9025
$ws.Range("I132").Value = 1068.5807
$ws.Range("K132").Value = 3205.7421
$ws.Range("M132").Value = -675.7420999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1666.1904
$ws.Range("I107").Value = 1070.2
$ws.Range("K107").Value = 3210.6
$ws.Range("M107").Value = -1290.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2457363.5
$ws.Range("I113").Value = 419
$ws.Range("J113").Value = 13513614
$ws.Range("K113").Value = 1257
$ws.Range("L113").Value = 40540842
$ws.Range("M113").Value = 913
$ws.Range("N113").Value = -40545182

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1366.9143
$ws.Range("I132").Value = 1204.8966
$ws.Range("J132").Value = 2150
$ws.Range("K132").Value = 3614.6898
$ws.Range("L132").Value = 6450
$ws.Range("M132").Value = -1084.6898
$ws.Range("N132").Value = -11510

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1071.2
$ws.Range("I136").Value = 621.6667
$ws.Range("J136").Value = 2419.8
$ws.Range("K136").Value = 1865.0001
$ws.Range("L136").Value = 7259.400000000001
$ws.Range("M136").Value = 684.9999
$ws.Range("N136").Value = -12359.4
